$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 4,5,6 (inserted before old row 7) ---
$ws.Range("B4").Value = 2024
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = 2026
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1

# --- New data row 8 (inserted before old row 9) ---
$ws.Range("B8").Value = 2009
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 400

# --- New data rows 10, 11 ---
$ws.Range("B10").Value = 2009
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = -1

$ws.Range("B11").Value = 2009
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 1

# --- New error-producing rows 12, 13, 14 ---
$ws.Range("B12").Value = 1900
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = -1

$ws.Range("B13").Value = 10000
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1

$ws.Range("B14").Value = -1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1

# --- New rows 15, 16 with quote-prefixed text "1" ---
$ws.Range("B15").Value = 2016
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'1"

$ws.Range("B16").Value = 2016
$ws.Range("C16").Value = "'1"
$ws.Range("D16").Value = 1

# --- Make sure every A3:A17 cell carries the same date number-format/style as A2 ---
# --- (done before writing formulas so the engine never has to manufacture a ---
# --- throwaway number-format style for the freshly computed date serials) ---
$ws.Range("A2").Copy()
$ws.Range("A3:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Apply the DATE formula across the full range A3:A17 in one shot so the ---
# --- engine groups identical formulas as shared formulas and recomputes all values. ---
$ws.Range("A3:A17").Formula = "=DATE(B3,C3,D3)"

# --- New block rows 22-24 ---
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 6

# --- Row 25: A25 standalone DATE formula, B25:F25 a shared DATE formula ---
$ws.Range("A2").Copy()
$ws.Range("A25:F25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A25").Formula = "=DATE(A24,A23,A22)"
$ws.Range("B25:F25").Formula = "=DATE(B24,B23,B22)"

# --- Selection matches the author's final cursor position ---
$ws.Range("C24").Select()

Write-Host "Done"
